$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New "find_simple_eratosfen_memopt" column (G) on the summary table
#    (rows 2-10). Numeric-looking values must be forced to text so they are
#    stored as shared strings (matching the source data), not as numbers.
# ---------------------------------------------------------------------------
function Set-TextValue($cellRef, $text) {
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range($cellRef).Style = "Normal"
}

Set-TextValue "G2" "find_simple_eratosfen_memopt"
Set-TextValue "G3" "0.0027307"
Set-TextValue "G4" "0.0028026"
Set-TextValue "G5" "0.0022225"
Set-TextValue "G6" "0.0031309"
Set-TextValue "G7" "0.0062043"
Set-TextValue "G8" "0.01344"
Set-TextValue "G9" "0.2914"
$ws.Range("G10").Value = 3.7217

# ---------------------------------------------------------------------------
# 2. New "Details" rows (60-67) for the find_simple_eratosfen_memopt scenario,
#    mirroring the existing per-testcase detail blocks above.
# ---------------------------------------------------------------------------
$detailRows = @(
    @{ Row=60; B="test.0.in"; C=0; E="0.0027307"; F="2020-02-07 15:02:05.803531" },
    @{ Row=61; B="test.1.in"; C=1; E="0.0028026"; F="2020-02-07 15:02:05.806037" },
    @{ Row=62; B="test.2.in"; C=2; E="0.0022225"; F="2020-02-07 15:02:05.808545" },
    @{ Row=63; B="test.3.in"; C=3; E="0.0031309"; F="2020-02-07 15:02:05.811553" },
    @{ Row=64; B="test.4.in"; C=4; E="0.0062043"; F="2020-02-07 15:02:05.818071" },
    @{ Row=65; B="test.5.in"; C=5; E="0.01344";   F="2020-02-07 15:02:05.831693" },
    @{ Row=66; B="test.6.in"; C=6; E="0.2914";    F="2020-02-07 15:02:06.122937" },
    @{ Row=67; B="test.7.in"; C=7; E=$null;       F="2020-02-07 15:02:09.845518" }
)

foreach ($r in $detailRows) {
    $row = $r.Row
    Set-TextValue "A$row" "find_simple_eratosfen_memopt"
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $true
    if ($r.E -ne $null) {
        Set-TextValue "E$row" $r.E
    } else {
        $ws.Range("E$row").Value = 3.7217
    }
    $ws.Range("F$row").Value = $r.F
}

# ---------------------------------------------------------------------------
# 3. Cosmetic sheet/workbook view updates.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 39.59
$ws.Columns.Item(7).ColumnWidth = 25.25

$excel.ActiveWindow.Zoom = 85
